# Updates rows 2-9 of Sheet1 with refreshed TPM-based NATMI ligand-receptor stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @("ECs", "a", "Mc4r", "ECs", 1, 0.3333333333333333, 0.045339, 0.136017, 0.1740293637846656, 0.1740293637846656, 1, 0.3333333333333333, 0.01244633333333333, 0.037339, 0.03341097222806328, 0.03341097222806328, 0.0005643043069999999, 0.005078738762999999, 0.005814490240276983, 0.005814490240276984),
    @("ECs", "a", "Mc4r", "FAPs", 1, 0.3333333333333333, 0.045339, 0.136017, 0.1740293637846656, 0.1740293637846656, 2, 0.6666666666666666, 0.05317933333333334, 0.159538, 0.1427547520640821, 0.1427547520640821, 0.002411097794, 0.021699880146, 0.02484351867894987, 0.02484351867894988),
    @("ECs", "a", "Mc4r", "MuSCs", 1, 0.3333333333333333, 0.045339, 0.136017, 0.1740293637846656, 0.1740293637846656, 3, 1, 0.282443, 0.847329, 0.758190784087218, 0.7581907840872181, 0.012805683177, 0.115251148593, 0.1319474597820953, 0.1319474597820953),
    @("ECs", "a", "Mc4r", "Resolving-Mac", 1, 0.3333333333333333, 0.045339, 0.136017, 0.1740293637846656, 0.1740293637846656, 1, 0.3333333333333333, 0.02445366666666667, 0.073361, 0.06564349162063661, 0.06564349162063661, 0.001108704793, 0.009978343136999999, 0.01142389508334342, 0.01142389508334342),
    @("MuSCs", "a", "Mc4r", "ECs", 1, 0.3333333333333333, 0.215186, 0.645558, 0.8259706362153344, 0.8259706362153345, 1, 0.3333333333333333, 0.01244633333333333, 0.037339, 0.03341097222806328, 0.03341097222806328, 0.002678276684666666, 0.024104490162, 0.0275964819877863, 0.0275964819877863),
    @("MuSCs", "a", "Mc4r", "FAPs", 1, 0.3333333333333333, 0.215186, 0.645558, 0.8259706362153344, 0.8259706362153345, 2, 0.6666666666666666, 0.05317933333333334, 0.159538, 0.1427547520640821, 0.1427547520640821, 0.01144344802266667, 0.102991032204, 0.1179112333851322, 0.1179112333851322),
    @("MuSCs", "a", "Mc4r", "MuSCs", 1, 0.3333333333333333, 0.215186, 0.645558, 0.8259706362153344, 0.8259706362153345, 3, 1, 0.282443, 0.847329, 0.758190784087218, 0.7581907840872181, 0.060777779398, 0.547000014582, 0.6262433243051226, 0.6262433243051229),
    @("MuSCs", "a", "Mc4r", "Resolving-Mac", 1, 0.3333333333333333, 0.215186, 0.645558, 0.8259706362153344, 0.8259706362153345, 1, 0.3333333333333333, 0.02445366666666667, 0.073361, 0.06564349162063661, 0.06564349162063661, 0.005262086715333333, 0.04735878043799999, 0.05421959653729319, 0.0542195965372932)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $rows[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $rowValues[$j]
    }
}

Write-Host "Updated rows 2-9 across columns A-T"
